# Added a poll and amended two old polls
# Column H ("CAN" = national numbers) on Sheet1 is updated:
#  - LIB/CON/NDP/BQ rows (4-7): new poll numbers, value + style normalized
#  - GRN/OTH rows (8-9): value unchanged, style normalized
#  - nw/nu rows (10-11): new sample-size totals

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PollValue($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.Value = $value
    # Re-stamp the (unchanged) font explicitly - this is what normalizes the
    # cell's style index the same way the source workbook's edit did.
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
}

Set-PollValue "H4" 35
Set-PollValue "H5" 30
Set-PollValue "H6" 20
Set-PollValue "H7" 7
Set-PollValue "H8" 6
Set-PollValue "H9" 2

$ws.Range("H10").Value = 1238
$ws.Range("H11").Value = 1242

$ws.Range("I4").Select() | Out-Null
